{"js": "const replacements = [\n  [\"2024-12-14 Saturday\", \"2024-12-15 Sunday\"],\n  [\"47\u00d782=\", \"15\u00d727=\"],\n  [\"28\u00d727=\", \"13\u00d740=\"],\n  [\"46\u00d747=\", \"72\u00d725=\"],\n  [\"58\u00d767=\", \"31\u00d753=\"],\n  [\"98\u00d791=\", \"59\u00d781=\"],\n  [\"79\u00d756=\", \"78\u00d766=\"],\n  [\"53\u00d715=\", \"15\u00d794=\"],\n  [\"38\u00d713=\", \"81\u00d789=\"],\n  [\"39\u00d726=\", \"54\u00d732=\"],\n  [\"11\u00d725=\", \"21\u00d795=\"],\n  [\"76\u00d713=\", \"95\u00d768=\"],\n  [\"71\u00d766=\", \"87\u00d717=\"],\n  [\"99\u00d760=\", \"34\u00d780=\"],\n  [\"69\u00d761=\", \"98\u00d716=\"],\n  [\"29\u00d759=\", \"16\u00d752=\"],\n  [\"70\u00d788=\", \"54\u00d731=\"],\n  [\"72\u00d740=\", \"24\u00d732=\"],\n  [\"61\u00d746=\", \"26\u00d729=\"],\n  [\"57\u00d750=\", \"36\u00d774=\"],\n  [\"62\u00d729=\", \"54\u00d756=\"],\n  [\"89\u00d793=\", \"48\u00d717=\"],\n  [\"69\u00d721=\", \"63\u00d765=\"],\n  [\"25\u00d737=\", \"28\u00d712=\"],\n  [\"40\u00d736=\", \"59\u00d777=\"],\n  [\"35\u00d746=\", \"99\u00d747=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"2024-12-14 Saturday\"; New = \"2024-12-15 Sunday\" },\n    @{ Old = \"47\u00d782=\"; New = \"15\u00d727=\" },\n    @{ Old = \"28\u00d727=\"; New = \"13\u00d740=\" },\n    @{ Old = \"46\u00d747=\"; New = \"72\u00d725=\" },\n    @{ Old = \"58\u00d767=\"; New = \"31\u00d753=\" },\n    @{ Old = \"98\u00d791=\"; New = \"59\u00d781=\" },\n    @{ Old = \"79\u00d756=\"; New = \"78\u00d766=\" },\n    @{ Old = \"53\u00d715=\"; New = \"15\u00d794=\" },\n    @{ Old = \"38\u00d713=\"; New = \"81\u00d789=\" },\n    @{ Old = \"39\u00d726=\"; New = \"54\u00d732=\" },\n    @{ Old = \"11\u00d725=\"; New = \"21\u00d795=\" },\n    @{ Old = \"76\u00d713=\"; New = \"95\u00d768=\" },\n    @{ Old = \"71\u00d766=\"; New = \"87\u00d717=\" },\n    @{ Old = \"99\u00d760=\"; New = \"34\u00d780=\" },\n    @{ Old = \"69\u00d761=\"; New = \"98\u00d716=\" },\n    @{ Old = \"29\u00d759=\"; New = \"16\u00d752=\" },\n    @{ Old = \"70\u00d788=\"; New = \"54\u00d731=\" },\n    @{ Old = \"72\u00d740=\"; New = \"24\u00d732=\" },\n    @{ Old = \"61\u00d746=\"; New = \"26\u00d729=\" },\n    @{ Old = \"57\u00d750=\"; New = \"36\u00d774=\" },\n    @{ Old = \"62\u00d729=\"; New = \"54\u00d756=\" },\n    @{ Old = \"89\u00d793=\"; New = \"48\u00d717=\" },\n    @{ Old = \"69\u00d721=\"; New = \"63\u00d765=\" },\n    @{ Old = \"25\u00d737=\"; New = \"28\u00d712=\" },\n    @{ Old = \"40\u00d736=\"; New = \"59\u00d777=\" },\n    @{ Old = \"35\u00d746=\"; New = \"99\u00d747=\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $ok = $find.Execute($r.Old, $true, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)\n    if (-not $ok) {\n        throw \"Replacement failed for: $($r.Old)\"\n    }\n}\n"}
